# Update Model component diagram (docs/diagrams/ModelComponentClassDiagram.pptx)
#
# Rename the "UniquePersonList" box to "SortedUniquePersonList" (and widen it),
# and nudge the small decision-diamond / connector / callout shapes that sit
# around it so the diagram still reads cleanly.

function Get-ShapeById($slide, $id) {
    foreach ($sh in $slide.Shapes) {
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "UniquePersonList" rectangle -> "SortedUniquePersonList", widened ------
$rectPersonList = Get-ShapeById $s 49
$rectPersonList.Left = 348.0
$rectPersonList.Top = 224.20245361328125
$rectPersonList.Width = 121.24212646484375
$rectPersonList.Height = 27.303937911987305
$rectPersonList.TextFrame.TextRange.Text = "SortedUniquePersonList"

# --- Elbow Connector 29 (decision diamond 51 -> rectangle 49) --------------
$connector29 = Get-ShapeById $s 30
$connector29.Width = 15.688898086547852
$connector29.Height = 20.47795295715332

# --- Elbow Connector 58 (decision diamond 51 -> "UniqueTagList" 57) --------
$connector58 = Get-ShapeById $s 59
$connector58.Adjustments.Item(1) = 0.38139

# --- Small decision diamond next to "Person" rectangle ----------------------
$decisionPerson = Get-ShapeById $s 63
$decisionPerson.Left = 469.24212646484375
$decisionPerson.Top = 232.11749267578125

# --- Elbow Connector 63 (decision diamond 63 -> "Person" rectangle 62) -----
$connectorPerson = Get-ShapeById $s 64
$connectorPerson.Left = 487.8285827636719
$connectorPerson.Top = 238.69654846191406
$connectorPerson.Width = 9.311339378356934
$connectorPerson.Height = 0.2469291388988495
$connectorPerson.VerticalFlip = -1

# --- "1" callout textbox near the bottom of the List association ----------
$calloutTextBox = Get-ShapeById $s 54
$calloutTextBox.Left = 334.71112060546875
$calloutTextBox.Top = 240.74850463867188
